$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D/E columns so values such as '66.657.01' or
# '10.50' are kept as literal text (matching the source inlineStr cells)
# instead of being auto-coerced to numbers by the COM Value setter.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.657.01'
$ws.Range("E2").Value = '  +2.49%  '
$ws.Range("D3").Value = '3.732.50'
$ws.Range("E3").Value = '  +6.49%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '418.70'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '130.92'
$ws.Range("E6").Value = '  -0.98%  '
$ws.Range("D7").Value = '3.723.21'
$ws.Range("E7").Value = '  +6.52%  '
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").Value = '0.181'
$ws.Range("E11").Value = '  +11.62%  '
$ws.Range("D12").Value = '0.0000398'
$ws.Range("E12").Value = '  +53.10%  '
$ws.Range("D13").Value = '42.74'
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '10.50'
$ws.Range("E14").Value = '  +4.88%  '
$ws.Range("D15").Value = '4.317.03'
$ws.Range("E15").Value = '  +6.47%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '20.71'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '3.728.18'
$ws.Range("E18").Value = '  +5.88%  '
$ws.Range("D19").Value = '13.17'
$ws.Range("E19").Value = '  +4.25%  '
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").Value = '66.745.34'
$ws.Range("E21").Value = '  +2.76%  '
$ws.Range("D22").Value = '444.08'
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("D23").Value = '16.42'
$ws.Range("E23").Value = '  +23.58%  '
$ws.Range("D24").Value = '89.62'
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("D25").Value = '3.15'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("D26").Value = '38.55'
$ws.Range("E26").Value = '  +12.84%  '
$ws.Range("D27").Value = '10.20'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").Value = '3.33'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("D29").Value = '5.09'
$ws.Range("E29").Value = '  +4.66%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.124'
$ws.Range("E30").Value = '  +8.89%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '12.71'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D32").Value = '2.70'
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = '7.23'
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("D35").Value = '42.06'
$ws.Range("E35").Value = '  +6.13%  '
$ws.Range("D36").Value = '56.98'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '0.0491'
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("D39").Value = '0.0₃0742'
$ws.Range("E39").Value = '  +5.30%  '
$ws.Range("E40").Value = '  -3.83%  '
$ws.Range("D41").Value = '3.03'
$ws.Range("E41").Value = '  +30.43%  '
$ws.Range("D42").Value = '29.23'
$ws.Range("E42").Value = '  +34.79%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = '3.44'
$ws.Range("E44").Value = '  +4.12%  '
$ws.Range("D45").Value = '3.23'
$ws.Range("E45").Value = '  +30.90%  '
$ws.Range("D46").Value = '2.14'
$ws.Range("E46").Value = '  +6.40%  '
$ws.Range("D47").Value = '145.92'
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("E48").Value = '  -3.95%  '
$ws.Range("E49").Value = '  -5.70%  '
$ws.Range("D50").Value = '4.34'
$ws.Range("E50").Value = '  -4.21%  '
$ws.Range("D51").Value = '0.306'
$ws.Range("E51").Value = '  -2.46%  '

# Restore the default (unstyled) cell style now that the values are
# committed as text, so no stray number-format style sticks to the cells.
$textRange.Style = "Normal"
